$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy "Reporting Location" = "On Campus" down column B for rows 2 through 10
# (mirrors column A's formatting, which uses a quote-prefixed text style)
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 2).Value = "On Campus"
    $ws.Cells.Item($r, 1).Copy()
    $ws.Cells.Item($r, 2).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# Update the active selection to B9, matching the saved view state
$ws.Range("B9").Select() | Out-Null
